$wb = $excel.ActiveWorkbook

# --- Sheet 1: perform_score ---
$ws1 = $wb.Worksheets.Item("perform_score")

$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0.6179412663638073
$ws1.Range("E3").Value = 0.2631158740604405

$ws1.Range("C4").Value = 0.4032894770784745
$ws1.Range("D4").Value = 1
$ws1.Range("E4").Value = 0.09385465034257801

$ws1.Range("C7").Value = 0.5258542493189475
$ws1.Range("D7").Value = 0.7592750916259005
$ws1.Range("E7").Value = 0.1539538788908697

$ws1.Range("C8").Value = 0.5093697533141097
$ws1.Range("D8").Value = 0.7548668673553367
$ws1.Range("E8").Value = 0.1677684496388044

$ws1.Range("C9").Value = 0.5739394718338696
$ws1.Range("D9").Value = 0.7325124145915776
$ws1.Range("E9").Value = 0.1760031839503879

$ws1.Range("C10").Value = 0.4350867014405913
$ws1.Range("D10").Value = 0.8691106329719475
$ws1.Range("E10").Value = 0.1124507382518814

$ws1.Range("C11").Value = 0.5537265071792478
$ws1.Range("D11").Value = 0.8509751868720293
$ws1.Range("E11").Value = 0.09081965925284612

$ws1.Range("C12").Value = 0.6584834330021733
$ws1.Range("D12").Value = 0.5782207208179736
$ws1.Range("E12").Value = 0.2600230058648094

# --- Sheet 2: ranking ---
$ws2 = $wb.Worksheets.Item("ranking")

$ws2.Range("C3").Value = 3
$ws2.Range("D3").Value = 2

$ws2.Range("C4").Value = 2
$ws2.Range("D4").Value = 3
$ws2.Range("E4").Value = 1

$ws2.Range("C7").Value = 2
$ws2.Range("D7").Value = 3

$ws2.Range("C10").Value = 2
$ws2.Range("D10").Value = 3
$ws2.Range("E10").Value = 1
